$d = $word.ActiveDocument

# Locate the paragraph that ends with "Xin chào tất cả mọi người" and split
# a new paragraph right after that sentence (i.e. before the trailing
# bookmark that marks the last edit position).
$anchor = "Xin chào tất cả mọi người"
$findRng = $d.Content
$findRng.Find.Execute($anchor, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$findRng.Collapse(0)          # wdCollapseEnd: collapse to just after the matched text
$findRng.InsertBefore("`r")   # split the paragraph here, pushing the bookmark to the new paragraph

# The freshly created paragraph is now the last paragraph in the document.
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newText = "Chúc cô và các bạn có một ngày tốt lành"
$insertRng = $d.Range($newPara.Range.Start, $newPara.Range.Start)
$insertRng.InsertBefore($newText)

# Give the newly typed run the same vi-VN language formatting used elsewhere
# in the document (matches the other runs' <w:rPr><w:lang w:val="vi-VN"/></w:rPr>).
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$formatRng = $newPara.Range.Duplicate
$formatRng.MoveEnd(1, -1) | Out-Null   # exclude the trailing paragraph mark
$formatRng.LanguageID = "vi-VN"
